$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1. Delete the paragraph containing the italic Chinese phrase
#    "得兒女的名分" (the whole <w:p> goes away).
# --------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*得兒女的名分*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# --------------------------------------------------------------------
# 2. Delete the paragraph "This PDF version is provided under the same
#    license." entirely.
# --------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*This PDF version is provided under the same license.*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# --------------------------------------------------------------------
# 3. Rewrite the attribution paragraph that used to read:
#    "關鍵詞 (Biblica) (Chinese (Traditional)) is based on: Biblica Bible
#    Dictionary, Biblica, Inc., 2023, which is licensed under a CC BY-SA
#    4.0 license."
#    into the new Biblica Study Notes (Key Terms) copyright/adaptation
#    text, keeping the leading/trailing empty runs (and their rPr) in
#    place.
# --------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*is based on*Biblica Bible Dictionary*") {
        $target = $p
        break
    }
}

$r = $target.Range
$contentEnd = $r.End - 1
if ($contentEnd -lt $r.Start) {
    $contentEnd = $r.Start
}
$body = $d.Range($r.Start, $contentEnd)
$body.Delete()

$insertPoint = $d.Range($target.Range.Start, $target.Range.Start)
$insertPoint.InsertAfter("Biblica Study Notes (Key Terms)")
$insertPoint.Font.Bold = 1

$insertPoint = $d.Range($insertPoint.End, $insertPoint.End)
$insertPoint.InsertAfter(" © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. ")
$insertPoint.Font.Bold = 0

$insertPoint = $d.Range($insertPoint.End, $insertPoint.End)
$insertPoint.InsertAfter("Biblica Study Notes")
$insertPoint.Font.Bold = 0

$insertPoint = $d.Range($insertPoint.End, $insertPoint.End)
$insertPoint.InsertAfter(" has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual.")
$insertPoint.Font.Bold = 0

# --------------------------------------------------------------------
# 4. Delete the "License Information" Heading2 paragraph entirely.
# --------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*License Information*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}
